$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1894021.5
$ws.Range("I33").Value = 2272815
$ws.Range("K33").Value = 2272815
$ws.Range("M33").Value = -2272586

$ws.Range("H101").Value = 1168
$ws.Range("I101").Value = 1168
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 3504
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -1882

$ws.Range("H112").Value = 1271.4777
$ws.Range("J112").Value = 1271.4777
$ws.Range("L112").Value = 3814.4331
$ws.Range("N112").Value = -6030.4331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 66

$ws.Range("H9").Value = 27767
$ws.Range("J9").Value = 27767
$ws.Range("L9").Value = 27767
$ws.Range("N9").Value = -28107

$ws.Range("H20").Value = 27767
$ws.Range("J20").Value = 27767
$ws.Range("L20").Value = 27767
$ws.Range("N20").Value = -28307

$ws.Range("H23").Value = 25251
$ws.Range("J23").Value = 25251
$ws.Range("L23").Value = 25251
$ws.Range("N23").Value = -25769

$ws.Range("H32").Value = 4961.705
$ws.Range("I32").Value = 5735.45
$ws.Range("K32").Value = 5735.45
$ws.Range("M32").Value = -5448.45

$ws.Range("H37").Value = 32805.184
$ws.Range("I37").Value = 33406.8
$ws.Range("J37").Value = 32303.834
$ws.Range("K37").Value = 33406.8
$ws.Range("L37").Value = 32303.834
$ws.Range("M37").Value = -33133.8
$ws.Range("N37").Value = -32849.834

$ws.Range("H44").Value = 38495.332
$ws.Range("J44").Value = 38495.332
$ws.Range("L44").Value = 38495.332
$ws.Range("N44").Value = -39471.332

$ws.Range("H55").Value = 26261.666
$ws.Range("J55").Value = 26261.666
$ws.Range("L55").Value = 26261.666
$ws.Range("N55").Value = -26891.666

$ws.Range("H80").Value = 35425.332
$ws.Range("J80").Value = 35425.332
$ws.Range("L80").Value = 35425.332
$ws.Range("N80").Value = -37421.332

$ws.Range("H83").Value = 35425.332
$ws.Range("J83").Value = 35425.332
$ws.Range("L83").Value = 106275.996
$ws.Range("N83").Value = -116259.996

$ws.Range("H137").Value = 48719.8
$ws.Range("J137").Value = 48719.8
$ws.Range("L137").Value = 48719.8
$ws.Range("N137").Value = -58919.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 361.25
$ws.Range("I80").Value = 467.2
$ws.Range("K80").Value = 467.2
$ws.Range("M80").Value = 530.8

$ws.Range("H83").Value = 361.25
$ws.Range("I83").Value = 467.2
$ws.Range("K83").Value = 2336
$ws.Range("M83").Value = 2656

$ws.Range("H130").Value = 41867.5
$ws.Range("J130").Value = 41867.5
$ws.Range("L130").Value = 41867.5
$ws.Range("N130").Value = -51907.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12504457
$ws.Range("J99").Value = 6600
$ws.Range("L99").Value = 6600
$ws.Range("N99").Value = -9596

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws.Range("H126").Value = 12504457
$ws.Range("J126").Value = 6600
$ws.Range("L126").Value = 19800
$ws.Range("N126").Value = -24740

$ws.Range("H137").Value = 44390
$ws.Range("J137").Value = 44390
$ws.Range("L137").Value = 44390
$ws.Range("N137").Value = -54590

$ws.Range("H141").Value = 25445.455
$ws.Range("J141").Value = 25445.455
$ws.Range("L141").Value = 25445.455
$ws.Range("N141").Value = -35805.455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 810.3
$ws.Range("J131").Value = 843.3626
$ws.Range("L131").Value = 2530.0878
$ws.Range("N131").Value = -12610.0878

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H74").Value = 21237.5
$ws.Range("J74").Value = 21237.5
$ws.Range("L74").Value = 21237.5
$ws.Range("N74").Value = -23109.5

$ws.Range("H77").Value = 21237.5
$ws.Range("J77").Value = 21237.5
$ws.Range("L77").Value = 63712.5
$ws.Range("N77").Value = -73072.5

$ws.Range("H97").Value = 781
$ws.Range("I97").Value = 715.2857
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 715.2857
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -219.2857
$ws.Range("N97").Value = -2003

$ws.Range("H98").Value = 19643
$ws.Range("J98").Value = 19643
$ws.Range("L98").Value = 19643
$ws.Range("N98").Value = -25633

$ws.Range("H104").Value = 33000
$ws.Range("J104").Value = 33000
$ws.Range("L104").Value = 33000
$ws.Range("N104").Value = -39988

$ws.Range("H105").Value = 33000
$ws.Range("J105").Value = 33000
$ws.Range("L105").Value = 33000
$ws.Range("N105").Value = -39988

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2878.0557
$ws.Range("I7").Value = 2053.3333
$ws.Range("J7").Value = 7001.6665
$ws.Range("K7").Value = 2053.3333
$ws.Range("L7").Value = 7001.6665
$ws.Range("M7").Value = -1941.3333
$ws.Range("N7").Value = -7225.6665

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("N18").Value = 0

$ws.Range("H20").Value = 1966.6666
$ws.Range("I20").Value = 900
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = -674
$ws.Range("N20").Value = -2952

$ws.Range("H40").Value = 6278.857
$ws.Range("I40").Value = 5444.3335
$ws.Range("J40").Value = 7781
$ws.Range("K40").Value = 5444.3335
$ws.Range("L40").Value = 7781
$ws.Range("M40").Value = -5308.3335
$ws.Range("N40").Value = -8053

$ws.Range("H46").Value = 1504.5
$ws.Range("I46").Value = 718.2857
$ws.Range("J46").Value = 1927.8462
$ws.Range("K46").Value = 718.2857
$ws.Range("L46").Value = 1927.8462
$ws.Range("M46").Value = -530.2857
$ws.Range("N46").Value = -2303.8462

$ws.Range("H126").Value = 2878.0557
$ws.Range("I126").Value = 2053.3333
$ws.Range("J126").Value = 7001.6665
$ws.Range("K126").Value = 6159.999899999999
$ws.Range("L126").Value = 21004.9995
$ws.Range("M126").Value = -3689.999899999999
$ws.Range("N126").Value = -25944.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3905.3333
$ws.Range("I122").Value = 2022.8125
$ws.Range("J122").Value = 6056.7856
$ws.Range("K122").Value = 6068.4375
$ws.Range("L122").Value = 18170.3568
$ws.Range("M122").Value = -3618.4375
$ws.Range("N122").Value = -23070.3568

$ws.Range("H136").Value = 11526.883
$ws.Range("I136").Value = 11635.7
$ws.Range("K136").Value = 34907.10000000001
$ws.Range("M136").Value = -32357.10000000001
